$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the two "wool dryer balls" keyword rows (rows 2 and 3).
# This shifts the remaining keyword rows up by two, and since no other
# cell on the sheet references those shared strings any more, Excel will
# drop them from the shared-strings table on save.
$ws.Rows("2:3").Delete()

# Match the saved selection state from the authored workbook.
$ws.Range("H4").Select()
